# Auto-generated Excel COM-interop script
# Applies numeric cell updates to match the target diff across 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) of the FFXIV leve-profit workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2519.2727
$ws.Range("J19").Value = 2806.85
$ws.Range("L19").Value = 2806.85
$ws.Range("N19").Value = -3156.85
$ws.Range("H32").Value = 3665.7058
$ws.Range("I32").Value = 1747.4546
$ws.Range("J32").Value = 7182.5
$ws.Range("K32").Value = 1747.4546
$ws.Range("L32").Value = 7182.5
$ws.Range("M32").Value = -1421.4546
$ws.Range("N32").Value = -7834.5
$ws.Range("H98").Value = 1381.75
$ws.Range("I98").Value = 1024.44
$ws.Range("K98").Value = 1024.44
$ws.Range("M98").Value = 473.5599999999999
$ws.Range("H113").Value = 3953.1667
$ws.Range("I113").Value = 4318
$ws.Range("J113").Value = 3770.75
$ws.Range("K113").Value = 4318
$ws.Range("L113").Value = 3770.75
$ws.Range("M113").Value = -1064
$ws.Range("N113").Value = -10278.75
$ws.Range("H122").Value = 1381.75
$ws.Range("I122").Value = 1024.44
$ws.Range("K122").Value = 3073.32
$ws.Range("M122").Value = -623.3200000000002
$ws.Range("H132").Value = 2447.9124
$ws.Range("I132").Value = 2531.1482
$ws.Range("J132").Value = 949.6667
$ws.Range("K132").Value = 7593.444600000001
$ws.Range("L132").Value = 2849.0001
$ws.Range("M132").Value = -5063.444600000001
$ws.Range("N132").Value = -7909.0001
$ws.Range("H133").Value = 89570.14
$ws.Range("J133").Value = 89570.14
$ws.Range("L133").Value = 89570.14
$ws.Range("N133").Value = -99690.14

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11754.306
$ws.Range("I32").Value = 6721.5386
$ws.Range("K32").Value = 6721.5386
$ws.Range("M32").Value = -6434.5386
$ws.Range("H61").Value = 3654.4443
$ws.Range("I61").Value = 2844.7917
$ws.Range("J61").Value = 10131.667
$ws.Range("K61").Value = 2844.7917
$ws.Range("L61").Value = 10131.667
$ws.Range("M61").Value = -2632.7917
$ws.Range("N61").Value = -10555.667
$ws.Range("H63").Value = 3790.8
$ws.Range("I63").Value = 1488.5
$ws.Range("K63").Value = 1488.5
$ws.Range("M63").Value = -802.5
$ws.Range("H66").Value = 3790.8
$ws.Range("I66").Value = 1488.5
$ws.Range("K66").Value = 7442.5
$ws.Range("M66").Value = -4010.5
$ws.Range("H74").Value = 1989.1063
$ws.Range("I74").Value = 1801.9556
$ws.Range("K74").Value = 1801.9556
$ws.Range("M74").Value = -927.9556
$ws.Range("H77").Value = 1989.1063
$ws.Range("I77").Value = 1801.9556
$ws.Range("K77").Value = 9009.778
$ws.Range("M77").Value = -4641.778
$ws.Range("H122").Value = 7725.516
$ws.Range("I122").Value = 6990.875
$ws.Range("J122").Value = 10244.286
$ws.Range("K122").Value = 20972.625
$ws.Range("L122").Value = 30732.858
$ws.Range("M122").Value = -18522.625
$ws.Range("N122").Value = -35632.858
$ws.Range("H132").Value = 4553.5625
$ws.Range("I132").Value = 2928.7778
$ws.Range("J132").Value = 6642.5713
$ws.Range("K132").Value = 8786.3334
$ws.Range("L132").Value = 19927.7139
$ws.Range("M132").Value = -6256.3334
$ws.Range("N132").Value = -24987.7139
$ws.Range("H136").Value = 3654.4443
$ws.Range("I136").Value = 2844.7917
$ws.Range("J136").Value = 10131.667
$ws.Range("K136").Value = 8534.375100000001
$ws.Range("L136").Value = 30395.001
$ws.Range("M136").Value = -5984.375100000001
$ws.Range("N136").Value = -35495.001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3234.7778
$ws.Range("I20").Value = 3316.5
$ws.Range("J20").Value = 3169.4
$ws.Range("K20").Value = 3316.5
$ws.Range("L20").Value = 3169.4
$ws.Range("M20").Value = -3069.5
$ws.Range("N20").Value = -3663.4
$ws.Range("H22").Value = 599.5
$ws.Range("I22").Value = 599
$ws.Range("K22").Value = 599
$ws.Range("M22").Value = -426
$ws.Range("H107").Value = 2604.2083
$ws.Range("I107").Value = 2627.7
$ws.Range("K107").Value = 2627.7
$ws.Range("M107").Value = -707.6999999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4105.2373
$ws.Range("I31").Value = 2981.5208
$ws.Range("K31").Value = 2981.5208
$ws.Range("M31").Value = -2686.5208
$ws.Range("H34").Value = 4105.2373
$ws.Range("I34").Value = 2981.5208
$ws.Range("K34").Value = 2981.5208
$ws.Range("M34").Value = -2779.5208
$ws.Range("H99").Value = 7592.273
$ws.Range("I99").Value = 5936
$ws.Range("K99").Value = 5936
$ws.Range("M99").Value = -4438
$ws.Range("H122").Value = 3358.3
$ws.Range("I122").Value = 323
$ws.Range("J122").Value = 15499.5
$ws.Range("K122").Value = 969
$ws.Range("L122").Value = 46498.5
$ws.Range("M122").Value = 1481
$ws.Range("N122").Value = -51398.5
$ws.Range("H126").Value = 7592.273
$ws.Range("I126").Value = 5936
$ws.Range("K126").Value = 17808
$ws.Range("M126").Value = -15338
$ws.Range("H134").Value = 1413.75
$ws.Range("I134").Value = 1267.4375
$ws.Range("K134").Value = 3802.3125
$ws.Range("M134").Value = -1267.3125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 254
$ws.Range("J33").Value = 438
$ws.Range("L33").Value = 2628
$ws.Range("N33").Value = -3194
$ws.Range("H136").Value = 7064.3105
$ws.Range("I136").Value = 6468.8184
$ws.Range("J136").Value = 8935.857
$ws.Range("K136").Value = 19406.4552
$ws.Range("L136").Value = 26807.571
$ws.Range("M136").Value = -14306.4552
$ws.Range("N136").Value = -37007.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7083.1665
$ws.Range("I70").Value = 7639.8
$ws.Range("K70").Value = 7639.8
$ws.Range("M70").Value = -7369.8
$ws.Range("H73").Value = 7083.1665
$ws.Range("I73").Value = 7639.8
$ws.Range("K73").Value = 7639.8
$ws.Range("M73").Value = -6703.8
$ws.Range("H102").Value = 1264.6
$ws.Range("I102").Value = 1201.4828
$ws.Range("K102").Value = 1201.4828
$ws.Range("M102").Value = 420.5172
$ws.Range("H113").Value = 7798.8887
$ws.Range("I113").Value = 8148.75
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 8148.75
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -5978.75
$ws.Range("N113").Value = -9340
$ws.Range("H133").Value = 88497.85000000001
$ws.Range("J133").Value = 88497.85000000001
$ws.Range("L133").Value = 88497.85000000001
$ws.Range("N133").Value = -98617.85000000001
$ws.Range("H139").Value = 89832.75
$ws.Range("J139").Value = 89832.75
$ws.Range("L139").Value = 89832.75
$ws.Range("N139").Value = -100112.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 200
$ws.Range("K25").Value = 200
$ws.Range("M25").Value = 30
$ws.Range("H61").Value = 2614.5
$ws.Range("J61").Value = 2999.5
$ws.Range("L61").Value = 2999.5
$ws.Range("N61").Value = -3403.5
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459
$ws.Range("H113").Value = 2614.5
$ws.Range("J113").Value = 2999.5
$ws.Range("L113").Value = 2999.5
$ws.Range("N113").Value = -7339.5
$ws.Range("H132").Value = 6369.6514
$ws.Range("I132").Value = 8937.210999999999
$ws.Range("K132").Value = 26811.633
$ws.Range("M132").Value = -24281.633
$ws.Range("H136").Value = 4693.528
$ws.Range("I136").Value = 4096.1724
$ws.Range("K136").Value = 12288.5172
$ws.Range("M136").Value = -9738.517200000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 5592.143
$ws.Range("J74").Value = 5544.1665
$ws.Range("L74").Value = 5544.1665
$ws.Range("N74").Value = -7416.1665
$ws.Range("H77").Value = 5592.143
$ws.Range("J77").Value = 5544.1665
$ws.Range("L77").Value = 16632.4995
$ws.Range("N77").Value = -25992.4995
$ws.Range("H132").Value = 4081.56
$ws.Range("I132").Value = 3977.45
$ws.Range("J132").Value = 4498
$ws.Range("K132").Value = 11932.35
$ws.Range("L132").Value = 13494
$ws.Range("M132").Value = -9402.349999999999
$ws.Range("N132").Value = -18554
$ws.Range("H136").Value = 2513.0889
$ws.Range("I136").Value = 2329.976
$ws.Range("K136").Value = 6989.928
$ws.Range("M136").Value = -4439.928

Write-Output "Applied all cell updates."
